$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Assignments")

# --- Updates to existing rows (room re-assignments / capacity fixes) ---
# Each entry: row, new room (col B), optional new capacity (col E)

$ws.Range("B6").Value  = "B F1.25 Computer Lab"

$ws.Range("B27").Value = "B F1.2 - Class/ECON Lab"
$ws.Range("E27").Value = 20

$ws.Range("B28").Value = "B F2.27 Creative Writing and Translation Studio"
$ws.Range("E28").Value = 18

$ws.Range("B30").Value = "B F1.35 FBA Conference Room"

$ws.Range("B39").Value = "A B.2 - EE Lab"

$ws.Range("B52").Value = "B F1.23 - Amphitheater I"

$ws.Range("B58").Value = "B F1.23 - Amphitheater I"

$ws.Range("B60").Value = "B F1.23 - Amphitheater I"

$ws.Range("B64").Value = "A F1.3 - Computer Lab"

$ws.Range("B67").Value = "RC.G4 - GBE IV"

$ws.Range("B84").Value = "B F2.27 Creative Writing and Translation Studio"
$ws.Range("E84").Value = 18

$ws.Range("B97").Value = "A F1.18 - Computer Lab"

$ws.Range("B98").Value = "B F1.2 - Class/ECON Lab"

$ws.Range("B105").Value = "RC1.3 - GSM and Network Laboratories"

$ws.Range("B106").Value = "A F1.18 - Computer Lab"

$ws.Range("B107").Value = "A B.1 - VACD Multimedia Studio"

$ws.Range("B115").Value = "B F1.2 - Class/ECON Lab"
$ws.Range("E115").Value = 20

$ws.Range("B124").Value = "RC1.3 - GSM and Network Laboratories"

$ws.Range("B131").Value = "A F1.17"

$ws.Range("B132").Value = "A F1.18 - Computer Lab"

$ws.Range("B136").Value = "A F2.8 - Drawing Studio & A F2.16 - Architecture Studio"

$ws.Range("B144").Value = "B F2.17"

$ws.Range("B155").Value = "A F2.8 - Drawing Studio"

$ws.Range("B165").Value = "B F1.9"

$ws.Range("B166").Value = "A B.16 - VACD Drawing Studio"

$ws.Range("B172").Value = "A F2.8 - Drawing Studio & A F2.16 - Architecture Studio"

$ws.Range("B178").Value = "A F1.3 - Computer Lab"

$ws.Range("B184").Value = "B F1.10 Class/ART Studio"

$ws.Range("B187").Value = "B F1.35 FBA Conference Room"

$ws.Range("B188").Value = "B F2.16"

$ws.Range("B199").Value = "B F1.23 - Amphitheater I"

$ws.Range("B211").Value = "B F2.15 - Amphitheater II"

$ws.Range("B212").Value = "B F2.6"

$ws.Range("B219").Value = "B F1.22"

$ws.Range("B228").Value = "A B.8 - Fabrication Lab"

$ws.Range("B233").Value = "B F1.16"

$ws.Range("B246").Value = "B F2.27 Creative Writing and Translation Studio"
$ws.Range("E246").Value = 18

$ws.Range("B253").Value = "A B.13 - Class/PSY Lab"
$ws.Range("E253").Value = 12

$ws.Range("B258").Value = "B F1.22"

$ws.Range("B260").Value = "B F1.16"

$ws.Range("B264").Value = "A B.16 - VACD Drawing Studio"

$ws.Range("B285").Value = "B F2.16"

# --- Append new rows (graduate courses) starting at row 294 ---

$newRows = @(
    @("ARCH510.1", "A B.13 - Class/PSY Lab", "Wed. 17:00-18:50", 9, 12),
    @("ARCH517.1", "B F2.27 Creative Writing and Translation Studio", "Tue. 10:00-15:50", 15, 18),
    @("ARCH569.1", "A B.13 - Class/PSY Lab", "Thu. 17:00-18:50", 5, 12),
    @("ARCH570.1", "A B.13 - Class/PSY Lab", "Tue. 17:00-18:50", 5, 12),
    @("BIO513.1", "A F2.16 - Architecture Studio", "Fri. 14:00-16:50", 9, 20),
    @("BIO514.1", "B F2.27 Creative Writing and Translation Studio", "Wed. 17:00-19:50", 6, 18),
    @("BIO518.1", "A F2.16 - Architecture Studio", "Thu. 17:00-19:50", 9, 20),
    @("BIO604.1", "A F2.16 - Architecture Studio", "Wed. 17:00-19:50", 3, 20),
    @("BIO646.1", "A F2.16 - Architecture Studio", "Tue. 17:00-19:50", 3, 20),
    @("BUS602.1", "B F2.27 Creative Writing and Translation Studio", "Fri. 14:00-16:50", 2, 18),
    @("CS509.1", "B F2.27 Creative Writing and Translation Studio", "Tue. 17:00-19:50", 13, 18),
    @("CS511.1", "RC1.3 - GSM and Network Laboratories", "Thu. 17:00-19:50", 15, 20),
    @("ECON506.1", "A F3.8 - Big Architecture Studio", "Thu. 17:00-19:50", 2, 25),
    @("ECON601.1", "A F1.3 - Computer Lab", "Thu. 17:00-19:50", 2, 25),
    @("EDU583.1", "A B.13 - Class/PSY Lab", "Thu. 16:00-18:50", 7, 12),
    @("ELT562.1", "A F2.16 - Architecture Studio", "Fri. 17:00-19:50", 9, 20),
    @("ELT565.1", "B F1.2 - Class/ECON Lab", "Thu. 17:00-19:50", 9, 20),
    @("ELT599.1", "A B.13 - Class/PSY Lab", "Mon. 17:00-17:50", 4, 12),
    @("ELT660.1", "B F2.27 Creative Writing and Translation Studio", "Fri. 17:00-19:50", 1, 18),
    @("ELT670.1", "A B.13 - Class/PSY Lab", "Thu. 17:00-19:50", 1, 12),
    @("IBF507.1", "A B.2 - EE Lab", "Wed. 17:00-19:50", 22, 25),
    @("IBF562.1", "A F1.4 - Class/Laboratory", "Tue. 17:00-19:50", 27, 30),
    @("IE502.1", "A B.13 - Class/PSY Lab", "Fri. 17:00-19:50", 4, 12),
    @("IR520.1", "B F1.2 - Class/ECON Lab", "Tue. 17:00-19:50", 3, 20),
    @("IR651.1", "A F2.16 - Architecture Studio", "Mon. 17:00-19:50", 2, 20),
    @("IR652.1", "RC1.4 - Computer Laboratory", "Tue. 17:00-19:50", 1, 20),
    @("LAW530.1", "A B.13 - Class/PSY Lab", "Fri. 18:00-20:50", 4, 12),
    @("MBA525.1", "A B.2 - EE Lab", "Fri. 17:00-19:50", 21, 25),
    @("MBA535.1", "A B.1 - VACD Multimedia Studio", "Mon. 17:00-19:50", 31, 35),
    @("MBA581.1", "RC1.3 - GSM and Network Laboratories", "Wed. 17:00-19:50", 6, 20),
    @("ME510.1", "A F1.3 - Computer Lab", "Tue. 17:00-19:50", 3, 25),
    @("PSY519.1", "RC1.3 - GSM and Network Laboratories", "Tue. 17:00-19:50", 7, 20),
    @("PSY524.1", "B F2.27 Creative Writing and Translation Studio", "Mon. 17:00-19:50", 4, 18),
    @("PSY529.1", "RC1.4 - Computer Laboratory", "Thu. 17:00-19:50", 6, 20),
    @("SOC503.1", "B F1.2 - Class/ECON Lab", "Wed. 17:00-19:50", 3, 20),
    @("SPS509.1", "A B.13 - Class/PSY Lab", "Mon. 17:00-19:50", 1, 12),
    @("SPS603.1", "A B.13 - Class/PSY Lab", "Wed. 17:00-19:50", 1, 12),
    @("VA502.1", "A B.13 - Class/PSY Lab", "Thu. 18:00-20:50", 7, 12),
    @("VA517.1", "A B.13 - Class/PSY Lab", "Tue. 17:00-19:50", 7, 12),
    @("VA519.1", "A B.13 - Class/PSY Lab", "Mon. 18:00-20:50", 7, 12)
)

$startRow = 294
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = "Assigned"
}
